# Apply "Last Version of data cleaning" edit:
# - Insert a new row "Bulgaria" at row 3 (pushing Denmark, Germany, Spain, France down by one row)
# - Remove the old "Hungary" row, which after the insert sits at row 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3; existing rows 3-11 shift down to 4-12
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with Bulgaria's data
$ws.Cells.Item(3, 1).Value = "Bulgaria"
$ws.Cells.Item(3, 2).Value = 83.09999999999999
$ws.Cells.Item(3, 3).Value = 79.8
$ws.Cells.Item(3, 4).Value = 83.09999999999999
$ws.Cells.Item(3, 5).Value = 81.8
$ws.Cells.Item(3, 6).Value = 85.3
$ws.Cells.Item(3, 7).Value = 90.89999999999999
$ws.Cells.Item(3, 8).Value = 91.09999999999999
$ws.Cells.Item(3, 9).Value = 95.60000000000001
$ws.Cells.Item(3, 10).Value = 94.8
$ws.Cells.Item(3, 11).Value = 99.30000000000001
$ws.Cells.Item(3, 12).Value = 98.90000000000001
$ws.Cells.Item(3, 13).Value = 101.4

# After the insert, the old "Hungary" row (formerly row 7) now lives at row 8.
# Delete it so the remaining rows (Netherlands, Poland, Finland, Sweden) shift back up
# to their original row numbers (8,9,10,11).
$ws.Rows.Item(8).Delete()
